# Applies the weekly Fruta/Hortaliza data update for the Cereza sheet
# (Vega Monumental Concepción). This corresponds to adding a new week of
# price observations (two additional "Primera"/"Segunda" rows) and
# shifting the historical weekly blocks down accordingly, per the
# authoritative diff of the workbook XML.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 57
$ws.Cells.Item(57, 4).Value = 44553
$ws.Cells.Item(57, 11).Value = 'Lapins'
$ws.Cells.Item(57, 13).Value = 250
$ws.Cells.Item(57, 14).Value = 5000
$ws.Cells.Item(57, 15).Value = 5500
$ws.Cells.Item(57, 16).Value = 5200
$ws.Cells.Item(57, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(57, 19).Value = 520

# Row 58
$ws.Cells.Item(58, 4).Value = 44553
$ws.Cells.Item(58, 11).Value = 'Rainier'
$ws.Cells.Item(58, 12).Value = 'Primera'
$ws.Cells.Item(58, 13).Value = 220
$ws.Cells.Item(58, 14).Value = 7000
$ws.Cells.Item(58, 15).Value = 7500
$ws.Cells.Item(58, 16).Value = 7227
$ws.Cells.Item(58, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(58, 19).Value = 723

# Row 59
$ws.Cells.Item(59, 4).Value = 44160
$ws.Cells.Item(59, 11).Value = 'Royal Dawn'
$ws.Cells.Item(59, 13).Value = 200
$ws.Cells.Item(59, 14).Value = 15000
$ws.Cells.Item(59, 15).Value = 16000
$ws.Cells.Item(59, 16).Value = 15500
$ws.Cells.Item(59, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(59, 19).Value = 1550

# Row 60
$ws.Cells.Item(60, 4).Value = 44160
$ws.Cells.Item(60, 11).Value = 'Royal Dawn'
$ws.Cells.Item(60, 13).Value = 100
$ws.Cells.Item(60, 14).Value = 13000
$ws.Cells.Item(60, 15).Value = 13000
$ws.Cells.Item(60, 16).Value = 13000
$ws.Cells.Item(60, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(60, 19).Value = 1300

# Row 61
$ws.Cells.Item(61, 4).Value = 44208
$ws.Cells.Item(61, 13).Value = 100
$ws.Cells.Item(61, 14).Value = 10000
$ws.Cells.Item(61, 15).Value = 11000
$ws.Cells.Item(61, 16).Value = 10500
$ws.Cells.Item(61, 18).Value = 'Región de Ñuble'
$ws.Cells.Item(61, 19).Value = 1050

# Row 62
$ws.Cells.Item(62, 4).Value = 44208
$ws.Cells.Item(62, 13).Value = 50
$ws.Cells.Item(62, 14).Value = 9000
$ws.Cells.Item(62, 15).Value = 9000
$ws.Cells.Item(62, 16).Value = 9000
$ws.Cells.Item(62, 18).Value = 'Región de Ñuble'
$ws.Cells.Item(62, 19).Value = 900

# Row 63
$ws.Cells.Item(63, 11).Value = 'Lapins'
$ws.Cells.Item(63, 13).Value = 200
$ws.Cells.Item(63, 14).Value = 9000
$ws.Cells.Item(63, 15).Value = 10000
$ws.Cells.Item(63, 16).Value = 9500
$ws.Cells.Item(63, 19).Value = 950

# Row 64
$ws.Cells.Item(64, 11).Value = 'Lapins'
$ws.Cells.Item(64, 13).Value = 100

# Row 65
$ws.Cells.Item(65, 4).Value = 44187
$ws.Cells.Item(65, 11).Value = 'Rainier'
$ws.Cells.Item(65, 13).Value = 100

# Row 66
$ws.Cells.Item(66, 4).Value = 44187
$ws.Cells.Item(66, 11).Value = 'Rainier'
$ws.Cells.Item(66, 13).Value = 50
$ws.Cells.Item(66, 14).Value = 8000
$ws.Cells.Item(66, 15).Value = 8000
$ws.Cells.Item(66, 16).Value = 8000
$ws.Cells.Item(66, 19).Value = 800

# Row 67
$ws.Cells.Item(67, 4).Value = 44196
$ws.Cells.Item(67, 13).Value = 200

# Row 68
$ws.Cells.Item(68, 4).Value = 44196
$ws.Cells.Item(68, 13).Value = 100

# Row 69
$ws.Cells.Item(69, 11).Value = 'Lapins'
$ws.Cells.Item(69, 13).Value = 100

# Row 70
$ws.Cells.Item(70, 11).Value = 'Lapins'
$ws.Cells.Item(70, 13).Value = 50

# Row 71
$ws.Cells.Item(71, 4).Value = 44188
$ws.Cells.Item(71, 11).Value = 'Rainier'
$ws.Cells.Item(71, 13).Value = 200

# Row 72
$ws.Cells.Item(72, 4).Value = 44188
$ws.Cells.Item(72, 11).Value = 'Rainier'
$ws.Cells.Item(72, 13).Value = 100

# Row 73
$ws.Cells.Item(73, 4).Value = 44195
$ws.Cells.Item(73, 13).Value = 100
$ws.Cells.Item(73, 14).Value = 10000
$ws.Cells.Item(73, 15).Value = 11000
$ws.Cells.Item(73, 16).Value = 10500
$ws.Cells.Item(73, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(73, 19).Value = 1050

# Row 74
$ws.Cells.Item(74, 4).Value = 44195
$ws.Cells.Item(74, 13).Value = 50
$ws.Cells.Item(74, 14).Value = 9000
$ws.Cells.Item(74, 15).Value = 9000
$ws.Cells.Item(74, 16).Value = 9000
$ws.Cells.Item(74, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(74, 19).Value = 900

# Row 75
$ws.Cells.Item(75, 11).Value = 'Lapins'
$ws.Cells.Item(75, 13).Value = 200
$ws.Cells.Item(75, 14).Value = 8000
$ws.Cells.Item(75, 15).Value = 8500
$ws.Cells.Item(75, 16).Value = 8250
$ws.Cells.Item(75, 19).Value = 825

# Row 76
$ws.Cells.Item(76, 11).Value = 'Lapins'
$ws.Cells.Item(76, 13).Value = 100
$ws.Cells.Item(76, 14).Value = 7500
$ws.Cells.Item(76, 15).Value = 7500
$ws.Cells.Item(76, 16).Value = 7500
$ws.Cells.Item(76, 19).Value = 750

# Row 77
$ws.Cells.Item(77, 4).Value = 44194
$ws.Cells.Item(77, 11).Value = 'Rainier'
$ws.Cells.Item(77, 13).Value = 100
$ws.Cells.Item(77, 14).Value = 9000
$ws.Cells.Item(77, 15).Value = 10000
$ws.Cells.Item(77, 16).Value = 9500
$ws.Cells.Item(77, 19).Value = 950

# Row 78
$ws.Cells.Item(78, 4).Value = 44194
$ws.Cells.Item(78, 11).Value = 'Rainier'
$ws.Cells.Item(78, 14).Value = 8000
$ws.Cells.Item(78, 15).Value = 8000
$ws.Cells.Item(78, 16).Value = 8000
$ws.Cells.Item(78, 19).Value = 800

# Row 79
$ws.Cells.Item(79, 4).Value = 44518
$ws.Cells.Item(79, 11).Value = 'Santina'
$ws.Cells.Item(79, 13).Value = 50
$ws.Cells.Item(79, 14).Value = 32000
$ws.Cells.Item(79, 15).Value = 32000
$ws.Cells.Item(79, 16).Value = 32000
$ws.Cells.Item(79, 17).Value = '$/caja 10 kilos'
$ws.Cells.Item(79, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(79, 19).Value = 3200

# Row 80
$ws.Cells.Item(80, 4).Value = 44518
$ws.Cells.Item(80, 12).Value = 'Segunda'
$ws.Cells.Item(80, 13).Value = 50
$ws.Cells.Item(80, 14).Value = 28000
$ws.Cells.Item(80, 15).Value = 28000
$ws.Cells.Item(80, 16).Value = 28000
$ws.Cells.Item(80, 17).Value = '$/caja 10 kilos'
$ws.Cells.Item(80, 19).Value = 2800

# New row 81 (inserted weekly entry)
$ws.Cells.Item(81, 1).Value = 11
$ws.Cells.Item(81, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(81, 3).Value = 'Bíobío'
$ws.Cells.Item(81, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(81, 4).Value = 44540
$ws.Cells.Item(81, 5).Value = 8
$ws.Cells.Item(81, 6).Value = 'Fruta'
$ws.Cells.Item(81, 7).Value = 100103
$ws.Cells.Item(81, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(81, 9).Value = 100103001
$ws.Cells.Item(81, 10).Value = 'Cereza'
$ws.Cells.Item(81, 11).Value = 'Lapins'
$ws.Cells.Item(81, 12).Value = 'Primera'
$ws.Cells.Item(81, 13).Value = 220
$ws.Cells.Item(81, 14).Value = 8500
$ws.Cells.Item(81, 15).Value = 9000
$ws.Cells.Item(81, 16).Value = 8727
$ws.Cells.Item(81, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(81, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(81, 19).Value = 873
$ws.Cells.Item(81, 20).Value = 10

# New row 82 (inserted weekly entry)
$ws.Cells.Item(82, 1).Value = 11
$ws.Cells.Item(82, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(82, 3).Value = 'Bíobío'
$ws.Cells.Item(82, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(82, 4).Value = 44540
$ws.Cells.Item(82, 5).Value = 8
$ws.Cells.Item(82, 6).Value = 'Fruta'
$ws.Cells.Item(82, 7).Value = 100103
$ws.Cells.Item(82, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(82, 9).Value = 100103001
$ws.Cells.Item(82, 10).Value = 'Cereza'
$ws.Cells.Item(82, 11).Value = 'Santina'
$ws.Cells.Item(82, 12).Value = 'Primera'
$ws.Cells.Item(82, 13).Value = 80
$ws.Cells.Item(82, 14).Value = 8000
$ws.Cells.Item(82, 15).Value = 9000
$ws.Cells.Item(82, 16).Value = 8625
$ws.Cells.Item(82, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(82, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(82, 19).Value = 862
$ws.Cells.Item(82, 20).Value = 10
